$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell with the new shared string
$ws.Range("E1").Value = "Matières enseignés"

# Set explicit column widths for C, D, E (matches new <cols> entries)
$ws.Columns.Item(3).ColumnWidth = 26.736979166666668
$ws.Columns.Item(4).ColumnWidth = 14.877604166666666
$ws.Columns.Item(5).ColumnWidth = 30.877604166666668

# Move/restore the active selection to E6
$ws.Range("E6").Select()
